$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated VCI3M forecast data for December (dates shifted, values refreshed)
$data = @{
    1 = @(45627,45634,45641,45648,45655,45662,45669,45676,45683,45690,45697)
    2 = @(70.59999999999999,69.5,68.59999999999999,67.90000000000001,67.5,67.3,67.2,67.2,67.3,67.2,66.90000000000001)
    3 = @(54.9,55.9,57.2,58.7,60.6,62.7,65.09999999999999,67.7,70.40000000000001,72.90000000000001,75.2)
    4 = @(51.9,51.1,50.3,49.7,49.3,49.2,49.2,49.4,49.6,50,50.3)
    5 = @(74.40000000000001,74.09999999999999,73.40000000000001,72.59999999999999,71.59999999999999,70.5,69.3,68,66.7,65.3,64)
    6 = @(63.1,61.1,60.4,61,62.8,65.90000000000001,70,74.7,79.8,84.8,89.40000000000001)
    7 = @(69.90000000000001,67.40000000000001,65.59999999999999,64.59999999999999,64.09999999999999,64.09999999999999,64.2,64.40000000000001,64.5,64.40000000000001,64)
    8 = @(80.8,80.3,79.8,79.2,78.7,78.2,77.8,77.3,76.7,75.90000000000001,74.90000000000001)
    9 = @(77.09999999999999,76,75,74,73,72,70.8,69.40000000000001,67.7,65.7,63.6)
    10 = @(79.90000000000001,79,77.59999999999999,75.8,73.7,71.40000000000001,69.09999999999999,66.7,64.5,62.3,60.5)
    11 = @(85.3,85,84.7,84.2,83.7,82.90000000000001,81.90000000000001,80.7,79.09999999999999,77.2,75.09999999999999)
    12 = @(72.3,72.40000000000001,72.09999999999999,71.59999999999999,70.90000000000001,70,69,67.90000000000001,66.7,65.40000000000001,64)
    13 = @(69.09999999999999,68.8,68.3,67.59999999999999,66.90000000000001,66.09999999999999,65.2,64.40000000000001,63.5,62.7,61.8)
    14 = @(77,73.3,70.3,68.2,67,66.59999999999999,66.8,67.40000000000001,68.2,68.90000000000001,69.5)
    15 = @(77.8,77.3,76.90000000000001,76.7,76.7,76.59999999999999,76.40000000000001,76,75.09999999999999,73.8,72.09999999999999)
    16 = @(63,62.4,61.7,61.1,60.5,60,59.7,59.5,59.5,59.6,59.8)
    17 = @(69.40000000000001,69,68.7,68.5,68.5,68.59999999999999,68.8,69.2,69.59999999999999,70,70.2)
    18 = @(74.3,74,73.5,73,72.3,71.5,70.5,69.3,67.90000000000001,66.3,64.59999999999999)
    19 = @(64.3,60.5,57.5,55.5,54.5,54.4,55.1,56.3,57.9,59.6,61.2)
    20 = @(72.40000000000001,70,68.59999999999999,68.3,68.7,69.7,71.09999999999999,72.40000000000001,73.40000000000001,73.90000000000001,73.8)
    21 = @(65.2,63.1,61.1,59.5,58.3,57.5,57.1,57,57.1,57.4,57.7)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
